$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.714.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").Value = "'3.016.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'550.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = "'134.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.87%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = "'3.012.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.15%  '
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = "'0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.17%  '
$ws.Range("D11").Value = "'6.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.13%  '
$ws.Range("D12").Value = "'0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.06%  '
$ws.Range("D13").Value = "'0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.14%  '
$ws.Range("D14").Value = "'34.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").Value = "'3.505.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").Value = "'61.809.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.77%  '
$ws.Range("E17").Value = '  -2.92%  '
$ws.Range("D18").Value = "'3.012.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.30%  '
$ws.Range("D19").Value = "'6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").Value = "'473.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.18%  '
$ws.Range("D21").Value = "'13.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.19%  '
$ws.Range("E22").Value = '  -4.91%  '
$ws.Range("D23").Value = "'7.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.04%  '
$ws.Range("D24").Value = "'80.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = "'7.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.34%  '
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").Value = "'25.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.42%  '
$ws.Range("E33").Value = '  -4.24%  '
$ws.Range("D34").Value = "'5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").Value = "'55.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.05%  '
$ws.Range("D36").Value = "'5.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.99%  '
$ws.Range("D37").Value = "'458.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.61%  '
$ws.Range("D38").Value = "'3.214.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").Value = "'0.0798"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").Value = "'0.0384"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.59%  '
$ws.Range("E41").Value = '  -2.19%  '
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("E43").Value = '  -11.04%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = "'25.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.98%  '
$ws.Range("E46").Value = '  -5.93%  '
$ws.Range("E47").Value = '  -4.58%  '
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("D49").Value = "'118.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.69%  '
$ws.Range("D50").Value = "'0.0₃0494"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.83%  '
$ws.Range("E51").Value = '  +6.55%  '
